# Applies the commit's changes to the deck:
#   1. Removes the "Shipping Order" slide (slide 18 - Rectangle shapes
#      "Shipping Order" / "Reefer Container" / "Cargo").
#   2. Refreshes the cached slide-number field text on the slide that
#      shifts from position 20 -> 19 as a result of the deletion.
#   3. Updates the cached date placeholder text ("9/14/19" -> "9/19/19")
#      that appears on the slide master, every slide layout, and the
#      notes master.

$p = $ppt.ActivePresentation

# --- 1. Delete the "Shipping Order" slide -------------------------------
$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -like "Shipping Order*") {
            $target = $i
        }
    }
}
if ($target -ne $null) {
    $p.Slides.Item($target).Delete()
}

# --- 2. Fix up the slide-number placeholder that is now slide 19 --------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $s = $p.Slides.Item($i)
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "20" -and $sh.Name -like "Slide Number*") {
            $sh.TextFrame.TextRange.Text = [string]$i
        }
    }
}

# --- 3. Update the cached date placeholder text on master/layouts/notes ---
function Update-DateShape($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $sh = $shapes.Item($k)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "9/14/19") {
            $sh.TextFrame.TextRange.Text = "9/19/19"
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShape $layouts.Item($li).Shapes
}

if ($p.HasNotesMaster) {
    Update-DateShape $p.NotesMaster.Shapes
}
